$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 22
$ws.Range("E2").Value = 7
$ws.Range("C3").Value = 'exploration'
$ws.Range("D3").Value = 12
$ws.Range("E3").Value = 1
$ws.Range("C4").Value = 'information generation'
$ws.Range("D4").Value = 11
$ws.Range("E4").Value = $null
$ws.Range("C5").Value = 'continuous operation'
$ws.Range("D5").Value = 11
$ws.Range("E5").Value = 1
$ws.Range("C6").Value = 'specialized tasks'
$ws.Range("D6").Value = 10
$ws.Range("C7").Value = 'versatility'
$ws.Range("D8").Value = 9
$ws.Range("E8").Value = 1
$ws.Range("C9").Value = 'basic needs'
$ws.Range("D9").Value = 8
$ws.Range("E9").Value = 24
$ws.Range("C10").Value = 'more possibilities'
$ws.Range("D10").Value = 8
$ws.Range("C11").Value = 'obstacle removal'
$ws.Range("D11").Value = 8
$ws.Range("E11").Value = 1
$ws.Range("C12").Value = 'efficiency'
$ws.Range("D12").Value = 7
$ws.Range("E12").Value = $null
$ws.Range("C13").Value = 'situational awareness'
$ws.Range("D13").Value = 7
$ws.Range("E13").Value = 1
$ws.Range("C14").Value = 'autonomy'
$ws.Range("D14").Value = 7
$ws.Range("C15").Value = 'readiness'
$ws.Range("D15").Value = 6
$ws.Range("E15").Value = $null
$ws.Range("C16").Value = 'speed'
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = $null
$ws.Range("C17").Value = 'sensor technology'
$ws.Range("E17").Value = $null
$ws.Range("C18").Value = 'size'
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 1
$ws.Range("C19").Value = 'limitations in handling complex or multiple tasks'
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 2
$ws.Range("C20").Value = 'remote controlled'
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 2
$ws.Range("C21").Value = 'data collection'
$ws.Range("D21").Value = 5
$ws.Range("E21").Value = $null
$ws.Range("C22").Value = 'multiplicability'
$ws.Range("E22").Value = $null
$ws.Range("C23").Value = 'flight'
$ws.Range("E23").Value = $null
$ws.Range("C24").Value = 'adaptability'
$ws.Range("D24").Value = 4
$ws.Range("E24").Value = 4
$ws.Range("C25").Value = 'inspection'
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = $null
$ws.Range("E26").Value = 1
$ws.Range("C27").Value = 'cost calculation'
$ws.Range("E27").Value = $null
$ws.Range("C28").Value = 'strength'
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = $null
$ws.Range("C29").Value = 'enhanced capabilities'
$ws.Range("E29").Value = $null
$ws.Range("C30").Value = 'live broadcasting'
$ws.Range("E30").Value = $null
$ws.Range("C31").Value = 'underwater rescue'
$ws.Range("E31").Value = 1
$ws.Range("C32").Value = 'mine clearance'
$ws.Range("E32").Value = 1
$ws.Range("C34").Value = 'assess damage'
$ws.Range("E34").Value = $null
$ws.Range("C35").Value = 'load capacity'
$ws.Range("D35").Value = 2
$ws.Range("D36").Value = 13
$ws.Range("E36").Value = 4
$ws.Range("C37").Value = 'faulty analyses'
$ws.Range("E37").Value = 1
$ws.Range("C38").Value = 'error risk'
$ws.Range("D38").Value = 9
$ws.Range("D39").Value = 7
$ws.Range("E39").Value = 2
$ws.Range("C40").Value = 'limited autonomy'
$ws.Range("E40").Value = 2
$ws.Range("C41").Value = 'specialization'
$ws.Range("D41").Value = 4
$ws.Range("E41").Value = $null
$ws.Range("C42").Value = 'immature technology'
$ws.Range("D42").Value = 3
$ws.Range("E42").Value = 1
$ws.Range("C43").Value = 'speed comparison'
$ws.Range("E43").Value = $null
$ws.Range("D44").Value = 2
$ws.Range("E44").Value = 16
$ws.Range("C45").Value = 'robot runtime'
$ws.Range("E45").Value = 2
$ws.Range("C46").Value = 'limited flexibility'
$ws.Range("E46").Value = 3
$ws.Range("C47").Value = 'risk'
$ws.Range("D47").Value = 1
$ws.Range("E47").Value = 2
$ws.Range("C48").Value = 'movement'
$ws.Range("D48").Value = 1
$ws.Range("E48").Value = $null
$ws.Range("A49").Value = 'SA'
$ws.Range("C49").Value = 'accessibility'
$ws.Range("D49").Value = 30
$ws.Range("E49").Value = 4
$ws.Range("C50").Value = 'operational efficiency'
$ws.Range("D50").Value = 30
$ws.Range("E50").Value = 2
$ws.Range("C51").Value = 'reliability'
$ws.Range("D51").Value = 23
$ws.Range("E51").Value = 6
$ws.Range("E53").Value = 6
$ws.Range("C54").Value = 'access'
$ws.Range("D54").Value = 20
$ws.Range("E54").Value = 6
$ws.Range("C55").Value = 'strength'
$ws.Range("D55").Value = 17
$ws.Range("E55").Value = 2
$ws.Range("C56").Value = 'search capabilities'
$ws.Range("D56").Value = 15
$ws.Range("E56").Value = 2
$ws.Range("C57").Value = 'replaceability'
$ws.Range("D57").Value = 15
$ws.Range("D58").Value = 14
$ws.Range("E58").Value = $null
$ws.Range("C59").Value = 'automation'
$ws.Range("D59").Value = 12
$ws.Range("E59").Value = 1
$ws.Range("C60").Value = 'focus on task'
$ws.Range("D60").Value = 10
$ws.Range("E60").Value = 3
$ws.Range("C61").Value = 'delivery of goods'
$ws.Range("D61").Value = 8
$ws.Range("E61").Value = 6
$ws.Range("A62").Value = 'R'
$ws.Range("C62").Value = 'potential physical harm'
$ws.Range("D62").Value = 11
$ws.Range("C63").Value = 'objective concerns'
$ws.Range("D63").Value = 8
$ws.Range("E63").Value = $null
$ws.Range("C64").Value = 'potential misuse'
$ws.Range("D64").Value = 8
$ws.Range("E64").Value = $null
$ws.Range("C65").Value = 'technical issues'
$ws.Range("E65").Value = 3
$ws.Range("D66").Value = 7
$ws.Range("C67").Value = 'misidentification'
$ws.Range("D67").Value = 5
$ws.Range("E67").Value = 1
$ws.Range("C68").Value = 'durability'
$ws.Range("D68").Value = 3
$ws.Range("E68").Value = 18
$ws.Range("C69").Value = 'unpredictability'
$ws.Range("D69").Value = 2
$ws.Range("C70").Value = 'safety'
$ws.Range("D70").Value = 1
$ws.Range("E70").Value = 1
$ws.Range("C71").Value = 'material'
$ws.Range("D71").Value = 1
$ws.Range("E71").Value = 9
$ws.Range("A72").Value = 'HRIP'
$ws.Range("C72").Value = 'collaborative support'
$ws.Range("D72").Value = 2
$ws.Range("E72").Value = 2
$ws.Range("A73").Value = 'HRIP'
$ws.Range("C73").Value = 'sustained performance'
$ws.Range("E73").Value = 1
$ws.Range("C74").Value = 'emotional resilience'
$ws.Range("D74").Value = 1
$ws.Range("E74").Value = 1
$ws.Range("A75").Value = 'HRIN'
$ws.Range("C75").Value = 'lack of emotional understanding'
$ws.Range("D75").Value = 7
$ws.Range("E75").Value = $null
$ws.Range("C76").Value = 'fear'
$ws.Range("D76").Value = 7
$ws.Range("C77").Value = 'reliance and trust'
$ws.Range("D77").Value = 5
$ws.Range("E77").Value = 1
$ws.Range("C78").Value = 'emotional coldness'
$ws.Range("D78").Value = 4
$ws.Range("A79").Value = 'AP'
$ws.Range("C79").Value = 'autonomous decisions'
$ws.Range("D79").Value = 5
$ws.Range("E79").Value = $null
$ws.Range("C80").Value = 'lack of emotionality'
$ws.Range("D80").Value = 4
$ws.Range("C81").Value = 'efficiency'
$ws.Range("D81").Value = 4
$ws.Range("A82").Value = 'AN'
$ws.Range("C82").Value = 'lack of empathy'
$ws.Range("D82").Value = 8
$ws.Range("C83").Value = 'errors in autonomy'
$ws.Range("D83").Value = 4
$ws.Range("E83").Value = 1
$ws.Range("D84").Value = 4
$ws.Range("C85").Value = 'perceived negative anthropomorphism'
$ws.Range("D85").Value = 3
$ws.Range("D87").Value = 1
